$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 2996
$ws.Range("J48").Value = 2996
$ws.Range("L48").Value = 8988
$ws.Range("N48").Value = -9572

$ws.Range("H56").Value = 2996
$ws.Range("J56").Value = 2996
$ws.Range("L56").Value = 8988
$ws.Range("N56").Value = -10056

$ws.Range("H64").Value = 6392.3335
$ws.Range("I64").Value = 6392.3335
$ws.Range("K64").Value = 6392.3335
$ws.Range("M64").Value = -6144.3335

$ws.Range("H67").Value = 6392.3335
$ws.Range("I67").Value = 6392.3335
$ws.Range("K67").Value = 6392.3335
$ws.Range("M67").Value = -5534.3335

$ws.Range("H116").Value = 4600.9473
$ws.Range("I116").Value = 4816.385
$ws.Range("K116").Value = 4816.385
$ws.Range("M116").Value = -1374.385

$ws.Range("H133").Value = 63680.375
$ws.Range("J133").Value = 63680.375
$ws.Range("L133").Value = 63680.375
$ws.Range("N133").Value = -73800.375

$ws.Range("H138").Value = 3308.8635
$ws.Range("I138").Value = 4756.375
$ws.Range("J138").Value = 2481.7144
$ws.Range("K138").Value = 14269.125
$ws.Range("L138").Value = 7445.1432
$ws.Range("M138").Value = -9129.125
$ws.Range("N138").Value = -17725.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 2307.5264
$ws.Range("I97").Value = 2016.2667
$ws.Range("K97").Value = 2016.2667
$ws.Range("M97").Value = -1520.2667

$ws.Range("H122").Value = 6514.2856
$ws.Range("I122").Value = 4300
$ws.Range("J122").Value = 7400
$ws.Range("K122").Value = 12900
$ws.Range("L122").Value = 22200
$ws.Range("M122").Value = -10450
$ws.Range("N122").Value = -27100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = ""
$ws.Range("N55").Value = 0

$ws.Range("H134").Value = 12926.934
$ws.Range("I134").Value = 6416.8096
$ws.Range("K134").Value = 19250.4288
$ws.Range("M134").Value = -16715.4288

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 465
$ws.Range("J7").Value = 550
$ws.Range("L7").Value = 1650
$ws.Range("N7").Value = -1874

$ws.Range("H44").Value = 749.44446
$ws.Range("I44").Value = 248.33333
$ws.Range("K44").Value = 744.99999
$ws.Range("M44").Value = -346.99999

$ws.Range("H68").Value = 1355.1111
$ws.Range("J68").Value = 1449.3334
$ws.Range("L68").Value = 4348.0002
$ws.Range("N68").Value = -5970.0002

$ws.Range("H71").Value = 1355.1111
$ws.Range("J71").Value = 1449.3334
$ws.Range("L71").Value = 13044.0006
$ws.Range("N71").Value = -21156.0006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 38791.273
$ws.Range("I20").Value = 5242.5
$ws.Range("J20").Value = 46246.555
$ws.Range("K20").Value = 5242.5
$ws.Range("L20").Value = 46246.555
$ws.Range("M20").Value = -4997.5
$ws.Range("N20").Value = -46736.555

$ws.Range("H24").Value = 35177.6
$ws.Range("J24").Value = 37975.11
$ws.Range("L24").Value = 37975.11
$ws.Range("N24").Value = -38321.11

$ws.Range("H46").Value = 5678.6665
$ws.Range("I46").Value = 1020.5
$ws.Range("J46").Value = 14995
$ws.Range("K46").Value = 1020.5
$ws.Range("L46").Value = 14995
$ws.Range("M46").Value = -864.5
$ws.Range("N46").Value = -15307

$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = ""
$ws.Range("N47").Value = 0

$ws.Range("H48").Value = 8000
$ws.Range("I48").Value = 6000
$ws.Range("J48").Value = 10000
$ws.Range("K48").Value = 6000
$ws.Range("L48").Value = 10000
$ws.Range("M48").Value = -5515
$ws.Range("N48").Value = -10970

$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").Value = ""

$ws.Range("H102").Value = 3621.037
$ws.Range("I102").Value = 3737.1738
$ws.Range("K102").Value = 3737.1738
$ws.Range("M102").Value = -2115.1738

$ws.Range("H122").Value = 1825.875
$ws.Range("I122").Value = 1531.6
$ws.Range("J122").Value = 2316.3333
$ws.Range("K122").Value = 4594.799999999999
$ws.Range("L122").Value = 6948.999899999999
$ws.Range("M122").Value = -2144.799999999999
$ws.Range("N122").Value = -11848.9999

$ws.Range("H126").Value = 6372.36
$ws.Range("I126").Value = 7531.1
$ws.Range("J126").Value = 5599.8667
$ws.Range("K126").Value = 22593.3
$ws.Range("L126").Value = 16799.6001
$ws.Range("M126").Value = -20123.3
$ws.Range("N126").Value = -21739.6001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 7502.5
$ws.Range("J3").Value = 7502.5
$ws.Range("L3").Value = 7502.5
$ws.Range("N3").Value = -7726.5

$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = ""
$ws.Range("N4").Value = 0

$ws.Range("H11").Value = 9769.385
$ws.Range("I11").Value = 10000
$ws.Range("J11").Value = 9625.25
$ws.Range("K11").Value = 10000
$ws.Range("L11").Value = 9625.25
$ws.Range("M11").Value = -9860
$ws.Range("N11").Value = -9905.25

$ws.Range("H15").Value = 7502.5
$ws.Range("J15").Value = 7502.5
$ws.Range("L15").Value = 7502.5
$ws.Range("N15").Value = -7842.5

$ws.Range("H16").Value = 5115.6
$ws.Range("I16").Value = 7060
$ws.Range("J16").Value = 2199
$ws.Range("K16").Value = 7060
$ws.Range("L16").Value = 2199
$ws.Range("M16").Value = -6890
$ws.Range("N16").Value = -2539

$ws.Range("H22").Value = 2868.6191
$ws.Range("I22").Value = 2462.6667
$ws.Range("K22").Value = 2462.6667
$ws.Range("M22").Value = -2167.6667

$ws.Range("H27").Value = 2868.6191
$ws.Range("I27").Value = 2462.6667
$ws.Range("K27").Value = 2462.6667
$ws.Range("M27").Value = -2355.6667

$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = ""
$ws.Range("N28").Value = 0

$ws.Range("H34").Value = 17304.8
$ws.Range("I34").Value = 1000
$ws.Range("J34").Value = 21381
$ws.Range("K34").Value = 1000
$ws.Range("L34").Value = 21381
$ws.Range("M34").Value = -828
$ws.Range("N34").Value = -21725

$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = ""
$ws.Range("N37").Value = 0

$ws.Range("H43").Value = 27728.5
$ws.Range("J43").Value = 27728.5
$ws.Range("L43").Value = 27728.5
$ws.Range("N43").Value = -28114.5

$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = ""
$ws.Range("N44").Value = 0

$ws.Range("H68").Value = 4767170.5
$ws.Range("I68").Value = 5000
$ws.Range("J68").Value = 5107325.5
$ws.Range("K68").Value = 5000
$ws.Range("L68").Value = 5107325.5
$ws.Range("M68").Value = -4251
$ws.Range("N68").Value = -5108823.5

$ws.Range("H71").Value = 4767170.5
$ws.Range("I71").Value = 5000
$ws.Range("J71").Value = 5107325.5
$ws.Range("K71").Value = 25000
$ws.Range("L71").Value = 25536627.5
$ws.Range("M71").Value = -21256
$ws.Range("N71").Value = -25544115.5

$ws.Range("H93").Value = 20090.182
$ws.Range("I93").Value = 21856.285
$ws.Range("J93").Value = 16999.5
$ws.Range("K93").Value = 21856.285
$ws.Range("L93").Value = 16999.5
$ws.Range("M93").Value = -20608.285
$ws.Range("N93").Value = -19495.5

$ws.Range("H122").Value = 6110.2593
$ws.Range("I122").Value = 5615.154
$ws.Range("J122").Value = 6570
$ws.Range("K122").Value = 16845.462
$ws.Range("L122").Value = 19710
$ws.Range("M122").Value = -14395.462
$ws.Range("N122").Value = -24610

$ws.Range("H125").Value = 129535.2
$ws.Range("J125").Value = 129535.2
$ws.Range("L125").Value = 129535.2
$ws.Range("N125").Value = -139375.2

$ws.Range("H132").Value = 1118515.8
$ws.Range("I132").Value = 4131.778
$ws.Range("K132").Value = 12395.334
$ws.Range("M132").Value = -9865.334000000001

$ws.Range("H136").Value = 10101.44
$ws.Range("I136").Value = 9919.357
$ws.Range("J136").Value = 10333.182
$ws.Range("K136").Value = 29758.071
$ws.Range("L136").Value = 30999.546
$ws.Range("M136").Value = -27208.071
$ws.Range("N136").Value = -36099.546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 14413.577
$ws.Range("I126").Value = 16818.637
$ws.Range("K126").Value = 50455.91099999999
$ws.Range("M126").Value = -47985.91099999999

$ws.Range("H132").Value = 8575.177
$ws.Range("I132").Value = 1667.3
$ws.Range("J132").Value = 18443.572
$ws.Range("K132").Value = 5001.9
$ws.Range("L132").Value = 55330.716
$ws.Range("M132").Value = -2471.9
$ws.Range("N132").Value = -60390.716

$ws.Range("H133").Value = 71494.664
$ws.Range("J133").Value = 71494.664
$ws.Range("L133").Value = 71494.664
$ws.Range("N133").Value = -81614.664

$ws.Range("H135").Value = 42499.5
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").Value = ""

$ws.Range("H136").Value = 13484.35
$ws.Range("I136").Value = 2223.0833
$ws.Range("K136").Value = 6669.249899999999
$ws.Range("M136").Value = -4119.249899999999

$ws.Range("H138").Value = 123331.664
$ws.Range("J138").Value = 169997.5
$ws.Range("N138").Value = -180277.5
